{"js": "// Map of original equation text -> updated equation text, as described by the diff.\nconst replacements = [\n  [\"512\u00d72=1024\", \"218\u00d75=1090\"],\n  [\"385\u00d73=1155\", \"995\u00d78=7960\"],\n  [\"869\u00d75=4345\", \"354\u00d72=708\"],\n  [\"388\u00d72=776\", \"762\u00d73=2286\"],\n  [\"702\u00d79=6318\", \"417\u00d77=2919\"],\n  [\"603\u00d73=1809\", \"669\u00d76=4014\"],\n  [\"187\u00d73=561\", \"294\u00d76=1764\"],\n  [\"390\u00d72=780\", \"375\u00d79=3375\"],\n  [\"714\u00d74=2856\", \"691\u00d74=2764\"],\n  [\"716\u00d77=5012\", \"169\u00d75=845\"],\n  [\"359\u00d72=718\", \"963\u00d79=8667\"],\n  [\"203\u00d77=1421\", \"829\u00d74=3316\"],\n  [\"898\u00d76=5388\", \"718\u00d75=3590\"],\n  [\"591\u00d78=4728\", \"381\u00d73=1143\"],\n  [\"871\u00d79=7839\", \"605\u00d76=3630\"],\n  [\"768\u00d72=1536\", \"185\u00d73=555\"],\n  [\"108\u00d78=864\", \"672\u00d78=5376\"],\n  [\"947\u00d74=3788\", \"438\u00d79=3942\"],\n  [\"574\u00d79=5166\", \"660\u00d76=3960\"],\n  [\"503\u00d79=4527\", \"251\u00d76=1506\"],\n  [\"875\u00d76=5250\", \"186\u00d78=1488\"],\n  [\"617\u00d72=1234\", \"379\u00d75=1895\"],\n  [\"985\u00d79=8865\", \"354\u00d73=1062\"],\n  [\"560\u00d78=4480\", \"912\u00d75=4560\"],\n  [\"640\u00d78=5120\", \"106\u00d79=954\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication equation in the table\n# with the new equation/result, as described by the diff.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"512\u00d72=1024\"; New = \"218\u00d75=1090\" },\n    @{ Old = \"385\u00d73=1155\"; New = \"995\u00d78=7960\" },\n    @{ Old = \"869\u00d75=4345\"; New = \"354\u00d72=708\" },\n    @{ Old = \"388\u00d72=776\";  New = \"762\u00d73=2286\" },\n    @{ Old = \"702\u00d79=6318\"; New = \"417\u00d77=2919\" },\n    @{ Old = \"603\u00d73=1809\"; New = \"669\u00d76=4014\" },\n    @{ Old = \"187\u00d73=561\";  New = \"294\u00d76=1764\" },\n    @{ Old = \"390\u00d72=780\";  New = \"375\u00d79=3375\" },\n    @{ Old = \"714\u00d74=2856\"; New = \"691\u00d74=2764\" },\n    @{ Old = \"716\u00d77=5012\"; New = \"169\u00d75=845\" },\n    @{ Old = \"359\u00d72=718\";  New = \"963\u00d79=8667\" },\n    @{ Old = \"203\u00d77=1421\"; New = \"829\u00d74=3316\" },\n    @{ Old = \"898\u00d76=5388\"; New = \"718\u00d75=3590\" },\n    @{ Old = \"591\u00d78=4728\"; New = \"381\u00d73=1143\" },\n    @{ Old = \"871\u00d79=7839\"; New = \"605\u00d76=3630\" },\n    @{ Old = \"768\u00d72=1536\"; New = \"185\u00d73=555\" },\n    @{ Old = \"108\u00d78=864\";  New = \"672\u00d78=5376\" },\n    @{ Old = \"947\u00d74=3788\"; New = \"438\u00d79=3942\" },\n    @{ Old = \"574\u00d79=5166\"; New = \"660\u00d76=3960\" },\n    @{ Old = \"503\u00d79=4527\"; New = \"251\u00d76=1506\" },\n    @{ Old = \"875\u00d76=5250\"; New = \"186\u00d78=1488\" },\n    @{ Old = \"617\u00d72=1234\"; New = \"379\u00d75=1895\" },\n    @{ Old = \"985\u00d79=8865\"; New = \"354\u00d73=1062\" },\n    @{ Old = \"560\u00d78=4480\"; New = \"912\u00d75=4560\" },\n    @{ Old = \"640\u00d78=5120\"; New = \"106\u00d79=954\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        Write-Host \"WARNING: could not find text '$($pair.Old)'\"\n    }\n}\n"}
